$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '252.12'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.68'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.550'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05689'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.464'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8067'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.042'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1429'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07321'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03152'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02937'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09268'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001669'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.222'
$ws.Range("D15").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0005811'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006457'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005055'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '18HotbitTokenHTBBestin24h'

$ws.Range("B20").Value = 'BitKan'

$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001051'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '19BitKanKAN'

$ws.Range("B21").Value = 'NitroEx'

$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0001501'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '20NitroExNTX'

$ws.Range("B22").Value = 'UpBots'

$ws.Range("C22").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0003200'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '21UpBotsUBXT'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.988'
$ws.Range("D23").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.111'
$ws.Range("D25").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1276'
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04132'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006914'
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = 'BKEXToken'

$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1048'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'

$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002972'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009566'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005642'
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.7853'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.01697'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01010'
$ws.Range("D50").Style = "Normal"
